# Update clear dev BOM
#  - mark the SN74LVC1G125DCKT alternate part with a note+link in I6
#  - add a new BOM line (row 27): Standoff / H1, with an alternate-part note+link
#  - leave the active selection on B27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- I6: flag the alternative part with a hyperlink, styled blue / no underline ---
$i6 = $ws.Range("I6")
$i6.Value2 = "alternative part"
$ws.Hyperlinks.Add($i6, "https://www.mouser.com/ProductDetail/Texas-Instruments/SN74LVC1G125DCKT", [Type]::Missing, [Type]::Missing, "alternative part")
$i6.Font.Underline = 0
$i6.Font.Color = 16711680

# --- Row 27: new BOM entry "Standoff" (designator H1) ---
$ws.Range("A27").Value2 = 26

$b27 = $ws.Range("B27")
$b27.HorizontalAlignment = 1
$b27.Value2 = "Standoff"
$ws.Hyperlinks.Add($b27, "https://www.mcmaster.com/standoffs/", [Type]::Missing, [Type]::Missing, "Standoff")
$b27.Font.Color = 16711680
$b27.Font.Underline = 0

$ws.Range("C27").Value2 = 1
$ws.Range("D27").Value2 = "H1"
$ws.Range("F27").Value2 = "N"

$i27 = $ws.Range("I27")
$i27.Value2 = "alternative part"
$ws.Hyperlinks.Add($i27, "https://www.mcmaster.com/standoffs/", [Type]::Missing, [Type]::Missing, "alternative part")
$i27.Font.Color = 16711680
$i27.Font.Underline = 0

# --- match the saved selection from the source edit ---
$ws.Range("B27").Select()
